$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from an existing "Date Solved" cell so the new
# cells reuse the same style (short date) instead of Excel creating a
# brand new number format entry.
$ws.Range("F106").Copy()
$ws.Range("F107:F108").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$solvedDate = (Get-Date -Year 2025 -Month 9 -Day 29).Date

# Row 107 - Reorder List
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = "Linked List"
$ws.Range("C107").Value = "Reorder List"
$ws.Range("D107").Value = "Medium"
$ws.Range("E107").Value = "Done"
$ws.Range("F107").Value = $solvedDate
$ws.Range("G107").Value = "O(n)"
$ws.Range("H107").Value = "O(1)"
$ws.Range("I107").Value = "Reverse & Merge"

# Row 108 - Flatten a multilevel doubly linked list
$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "Linked List"
$ws.Range("C108").Value = "Flatten a Doubly Linked List"
$ws.Range("D108").Value = "Medium"
$ws.Range("E108").Value = "Done"
$ws.Range("F108").Value = $solvedDate
$ws.Range("G108").Value = "O(n)"
$ws.Range("H108").Value = "O(1)"
$ws.Range("I108").Value = "DFS"

# Update selection to match new active cell
$null = $ws.Range("A109").Select()
